# Rename the existing total ("总计") sheet to "2022-Q1" (keeps its sheetId),
# then add a brand-new sheet right after it and name that one "总计"
# (matches the sheetId numbering seen in the target file: 2022-Q1 -> 5, 总计 -> 6).
$wb = $excel.ActiveWorkbook
$styleSrc = $wb.Worksheets.Item(2)
$totalWs = $wb.Worksheets.Item("总计")
$totalName = $totalWs.Name
$totalWs.Name = "2022-Q1"
$q1 = $totalWs
$newTotal = $wb.Worksheets.Add($null, $q1)
$newTotal.Name = $totalName

# Wipe the old 总计 rows/columns out of the renamed sheet -- it will be
# rebuilt from scratch with the fund-holdings table.
$q1.Cells.Clear()

# Grab a header-style cell (bold/centered/bordered, style index 2 in the
# original workbook) to stamp onto the new header row + index column.
$styleSrc.Range("B1").Copy()

$q1headers = @('基金代码', '基金名称', '基金规模', '股票总仓位', '仓位占比', '持有市值(亿元)', '仓位排名')
$col = 2
foreach ($h in $q1headers) {
    $c = $q1.Cells.Item(1, $col)
    $c.Value = $h
    $c.PasteSpecial(-4122)
    $col = $col + 1
}

# A2:H30 fund holdings data
$q1data = @(
    @(0, '005299', '万家成长优选灵活配置混合A', '24.43', '93.91', '2.36', '0.5765', 10),
    @(1, '010694', '万家内需增长一年持有期混合', '17.21', '94.85', '3.08', '0.5301', 9),
    @(2, '005300', '万家成长优选灵活配置混合C', '10.93', '93.91', '2.36', '0.2579', 10),
    @(3, '161039', '富国中证1000指数增强LOF', '21.72', '89.03', '0.86', '0.1868', 2),
    @(4, '002601', '中银证券价值精选灵活配置混合', '3.41', '93.74', '4.98', '0.1698', 4),
    @(5, '006132', '万家智造优势混合A', '4.82', '93.70', '2.69', '0.1297', 7),
    @(6, '002504', '鹏华金鼎灵活配置混合A', '2.49', '77.53', '3.58', '0.0891', 9),
    @(7, '007251', '广发睿享稳健增利混合', '3.69', '38.80', '1.65', '0.0609', 9),
    @(8, '011269', '中银证券优势制造股票型证券投资基金A', '1.39', '93.51', '4.36', '0.0606', 8),
    @(9, '006377', '广发趋势动力灵活配置混合', '0.82', '90.19', '6.30', '0.0517', 1),
    @(10, '002025', '广发聚盛灵活配置混合A', '7.09', '22.40', '0.66', '0.0468', 4),
    @(11, '005189', '海富通量化前锋股票A', '3.09', '89.94', '1.38', '0.0426', 5),
    @(12, '009766', '安信平稳双利3个月持有期混合A', '2.33', '39.45', '1.76', '0.0410', 8),
    @(13, '000433', '安信鑫发优选混合', '1.23', '67.20', '2.81', '0.0346', 3),
    @(14, '009500', '国寿安保高股息混合A', '0.96', '73.47', '3.06', '0.0294', 8),
    @(15, '005492', '农银汇理研究驱动灵活配置混合', '0.77', '62.76', '2.19', '0.0169', 10),
    @(16, '006133', '万家智造优势混合C', '0.52', '93.70', '2.69', '0.0140', 7),
    @(17, '005188', '海富通量化前锋股票C', '0.81', '89.94', '1.38', '0.0112', 5),
    @(18, '002505', '鹏华金鼎灵活配置混合C', '0.26', '77.53', '3.58', '0.0093', 9),
    @(19, '011270', '中银证券优势制造股票型证券投资基金C', '0.21', '93.51', '4.36', '0.0092', 8),
    @(20, '011149', '创金合信ESG责任投资股票A', '0.16', '87.53', '4.50', '0.0072', 3),
    @(21, '002026', '广发聚盛灵活配置混合C', '1.07', '22.40', '0.66', '0.0071', 4),
    @(22, '009767', '安信平稳双利3个月持有期混合C', '0.26', '39.45', '1.76', '0.0046', 8),
    @(23, '011150', '创金合信ESG责任投资股票C', '0.08', '87.53', '4.50', '0.0036', 3),
    @(24, '004913', '中银证券聚瑞混合A', '0.10', '32.71', '3.12', '0.0031', 5),
    @(25, '750005', '安信平稳增长混合A', '0.08', '65.16', '2.95', '0.0024', 8),
    @(26, '009501', '国寿安保高股息混合C', '0.03', '73.47', '3.06', '0.0009', 8),
    @(27, '004914', '中银证券聚瑞混合C', '0.02', '32.71', '3.12', '0.0006', 5),
    @(28, '002035', '安信平稳增长混合C', '0.00', '65.16', '2.95', $null, 8)
)

$r = 2
foreach ($row in $q1data) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 1).PasteSpecial(-4122)
    $q1.Cells.Item($r, 2).Value = "'" + $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = "'" + $row[3]
    $q1.Cells.Item($r, 5).Value = "'" + $row[4]
    $q1.Cells.Item($r, 6).Value = "'" + $row[5]
    if ($row[6] -eq $null) {
        $q1.Cells.Item($r, 7).Value = 0
    } else {
        $q1.Cells.Item($r, 7).Value = "'" + $row[6]
    }
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Rebuild the 总计 (grand total) sheet: header + 5 rows (adds the new 2022-Q1 row on top)
$totalHeaders = @('日期', '持有数量(只)', '持有市值(亿元)')
$col = 2
foreach ($h in $totalHeaders) {
    $c = $newTotal.Cells.Item(1, $col)
    $c.Value = $h
    $c.PasteSpecial(-4122)
    $col = $col + 1
}

$totalData = @(
    @(0, '2022-Q1', 29, 2.4),
    @(1, '2021-Q4', 4, 0.2),
    @(2, '2021-Q3', 6, 0.18),
    @(3, '2021-Q2', 3, 0.15),
    @(4, '2021-Q1', 2, 0.01)
)

$r = 2
foreach ($row in $totalData) {
    $newTotal.Cells.Item($r, 1).Value = $row[0]
    $newTotal.Cells.Item($r, 1).PasteSpecial(-4122)
    $newTotal.Cells.Item($r, 2).Value = $row[1]
    $newTotal.Cells.Item($r, 3).Value = $row[2]
    $newTotal.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$q1.Range("A1").Select()
$newTotal.Range("A1").Select()
